$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IPC PO (predicted output), DELTA and DELTA^2 columns with refreshed model results
$ws.Cells.Item(2, 3).Value = 30.40247058788219
$ws.Cells.Item(2, 4).Value = 0.6524705878821884
$ws.Cells.Item(2, 5).Value = 0.4257178680513285
$ws.Cells.Item(3, 3).Value = 29.84945777265597
$ws.Cells.Item(3, 4).Value = 0.009457772655970587
$ws.Cells.Item(3, 5).Value = 0.00008944946361202493
$ws.Cells.Item(4, 3).Value = 29.66142247796292
$ws.Cells.Item(4, 4).Value = -0.14857752203707975
$ws.Cells.Item(4, 5).Value = 0.022075280054678917
$ws.Cells.Item(5, 3).Value = 29.5982817534586
$ws.Cells.Item(5, 4).Value = -0.32171824654140124
$ws.Cells.Item(5, 5).Value = 0.10350263015767383
$ws.Cells.Item(6, 3).Value = 29.85312025586845
$ws.Cells.Item(6, 4).Value = -0.1268797441315499
$ws.Cells.Item(6, 5).Value = 0.016098469470887573
$ws.Cells.Item(7, 3).Value = 30.06600351530143
$ws.Cells.Item(7, 4).Value = 0.026003515301429303
$ws.Cells.Item(7, 5).Value = 0.0006761828080316679
$ws.Cells.Item(8, 3).Value = 30.34908901132387
$ws.Cells.Item(8, 4).Value = 0.1390890113238683
$ws.Cells.Item(8, 5).Value = 0.019345753071051165
$ws.Cells.Item(9, 3).Value = 30.52594039433468
$ws.Cells.Item(9, 4).Value = 0.3059403943346801
$ws.Cells.Item(9, 5).Value = 0.09359952488565958
$ws.Cells.Item(10, 3).Value = 30.35365911871458
$ws.Cells.Item(10, 4).Value = -0.02634088128541734
$ws.Cells.Item(10, 5).Value = 0.0006938420268924494
$ws.Cells.Item(11, 3).Value = 30.43898411740251
$ws.Cells.Item(11, 4).Value = -0.0010158825974926344
$ws.Cells.Item(11, 5).Value = 0.0000010320174518883818
$ws.Cells.Item(12, 3).Value = 30.55094216973905
$ws.Cells.Item(12, 4).Value = 0.07094216973904821
$ws.Cells.Item(12, 5).Value = 0.005032791447283928
$ws.Cells.Item(13, 3).Value = 30.68086014304607
$ws.Cells.Item(13, 4).Value = -0.009139856953932934
$ws.Cells.Item(13, 5).Value = 0.00008353698513835622
$ws.Cells.Item(14, 3).Value = 30.27902283158048
$ws.Cells.Item(14, 4).Value = -0.47097716841951964
$ws.Cells.Item(14, 5).Value = 0.22181949317246857
$ws.Cells.Item(15, 3).Value = 30.38232876759964
$ws.Cells.Item(15, 4).Value = -0.557671232400363
$ws.Cells.Item(15, 5).Value = 0.31099720344693965
$ws.Cells.Item(16, 3).Value = 30.69458768077862
$ws.Cells.Item(16, 4).Value = -0.25541231922137797
$ws.Cells.Item(16, 5).Value = 0.06523545281004309
$ws.Cells.Item(17, 3).Value = 31.06272285481336
$ws.Cells.Item(17, 4).Value = 0.042722854813359135
$ws.Cells.Item(17, 5).Value = 0.0018252423234033637
$ws.Cells.Item(18, 3).Value = 31.35281011065088
$ws.Cells.Item(18, 4).Value = 0.23281011065087753
$ws.Cells.Item(18, 5).Value = 0.05420054762127384
$ws.Cells.Item(19, 3).Value = 31.33141430812179
$ws.Cells.Item(19, 4).Value = 0.05141430812178882
$ws.Cells.Item(19, 5).Value = 0.0026434310796422398
$ws.Cells.Item(20, 3).Value = 31.6087785071903
$ws.Cells.Item(20, 4).Value = 0.22877850719030235
$ws.Cells.Item(20, 5).Value = 0.05233960535222323
$ws.Cells.Item(21, 3).Value = 31.92626103420487
$ws.Cells.Item(21, 4).Value = 0.346261034204872
$ws.Cells.Item(21, 5).Value = 0.11989670380862752
$ws.Cells.Item(22, 3).Value = 31.9100204118511
$ws.Cells.Item(22, 4).Value = 0.2600204118511016
$ws.Cells.Item(22, 5).Value = 0.0676106145792165
$ws.Cells.Item(23, 3).Value = 31.59474956815411
$ws.Cells.Item(23, 4).Value = -0.2852504318458884
$ws.Cells.Item(23, 5).Value = 0.08136780886826582
$ws.Cells.Item(24, 3).Value = 32.11077502296926
$ws.Cells.Item(24, 4).Value = -0.16922497703074413
$ws.Cells.Item(24, 5).Value = 0.028637092851055878
$ws.Cells.Item(25, 3).Value = 32.1607344515993
$ws.Cells.Item(25, 4).Value = -0.2892655484006994
$ws.Cells.Item(25, 5).Value = 0.08367455749155737
$ws.Cells.Item(26, 3).Value = 33.26761652402123
$ws.Cells.Item(26, 4).Value = 0.41761652402122706
$ws.Cells.Item(26, 5).Value = 0.17440356113557212
$ws.Cells.Item(27, 3).Value = 32.95359507319187
$ws.Cells.Item(27, 4).Value = 0.053595073191871734
$ws.Cells.Item(27, 5).Value = 0.002872431870442088
$ws.Cells.Item(28, 3).Value = 33.28701813607838
$ws.Cells.Item(28, 4).Value = 0.18701813607837892
$ws.Cells.Item(28, 5).Value = 0.034975783222231056
$ws.Cells.Item(29, 3).Value = 33.29908565264557
$ws.Cells.Item(29, 4).Value = -0.10091434735442562
$ws.Cells.Item(29, 5).Value = 0.010183705501969668
$ws.Cells.Item(30, 3).Value = 33.93100440488091
$ws.Cells.Item(30, 4).Value = 0.23100440488090612
$ws.Cells.Item(30, 5).Value = 0.0533630350743816
$ws.Cells.Item(31, 3).Value = 34.69771747294202
$ws.Cells.Item(31, 4).Value = 0.5977174729420156
$ws.Cells.Item(31, 5).Value = 0.3572661774601892
$ws.Cells.Item(32, 3).Value = 34.72062645199188
$ws.Cells.Item(32, 4).Value = 0.32062645199188466
$ws.Cells.Item(32, 5).Value = 0.10280132171690431
$ws.Cells.Item(33, 3).Value = 35.21850464189565
$ws.Cells.Item(33, 4).Value = 0.3185046418956503
$ws.Cells.Item(33, 5).Value = 0.10144520690907644
$ws.Cells.Item(34, 3).Value = 35.05499915420738
$ws.Cells.Item(34, 4).Value = -0.24500084579261738
$ws.Cells.Item(34, 5).Value = 0.06002541443909788
$ws.Cells.Item(35, 3).Value = 35.3473033739953
$ws.Cells.Item(35, 4).Value = -0.3526966260047004
$ws.Cells.Item(35, 5).Value = 0.12439490999509951
$ws.Cells.Item(36, 3).Value = 35.98922756969985
$ws.Cells.Item(36, 4).Value = -0.31077243030014756
$ws.Cells.Item(36, 5).Value = 0.09657950343466007
$ws.Cells.Item(37, 3).Value = 36.38170493970585
$ws.Cells.Item(37, 4).Value = -0.41829506029414887
$ws.Cells.Item(37, 5).Value = 0.17497075746648563
$ws.Cells.Item(38, 3).Value = 37.41691575368267
$ws.Cells.Item(38, 4).Value = 0.11691575368267593
$ws.Cells.Item(38, 5).Value = 0.01366929345918815
$ws.Cells.Item(39, 3).Value = 37.95126721373015
$ws.Cells.Item(39, 4).Value = 0.05126721373014931
$ws.Cells.Item(39, 5).Value = 0.0026283272036528096
$ws.Cells.Item(40, 3).Value = 38.5437647617419
$ws.Cells.Item(40, 4).Value = 0.04376476174189747
$ws.Cells.Item(40, 5).Value = 0.0019153543703250527
$ws.Cells.Item(41, 3).Value = 39.41892405364155
$ws.Cells.Item(41, 4).Value = 0.5189240536415483
$ws.Cells.Item(41, 5).Value = 0.26928217344777644
$ws.Cells.Item(42, 3).Value = 40.13086809550075
$ws.Cells.Item(42, 4).Value = 0.7308680955007532
$ws.Cells.Item(42, 5).Value = 0.534168173020898
$ws.Cells.Item(43, 3).Value = 40.42224859150753
$ws.Cells.Item(43, 4).Value = 0.5222485915075339
$ws.Cells.Item(43, 5).Value = 0.272743591331603
$ws.Cells.Item(44, 3).Value = 40.01027232520678
$ws.Cells.Item(44, 4).Value = -0.08972767479322385
$ws.Cells.Item(44, 5).Value = 0.00805105562379854
$ws.Cells.Item(45, 3).Value = 41.15194085867312
$ws.Cells.Item(45, 4).Value = 0.5519408586731203
$ws.Cells.Item(45, 5).Value = 0.3046387114728213
$ws.Cells.Item(46, 3).Value = 41.37871835482554
$ws.Cells.Item(46, 4).Value = 0.4787183548255385
$ws.Cells.Item(46, 5).Value = 0.22917126324687018
$ws.Cells.Item(47, 3).Value = 41.16526965467785
$ws.Cells.Item(47, 4).Value = -0.0347303453221528
$ws.Cells.Item(47, 5).Value = 0.0012061968861959808
$ws.Cells.Item(48, 3).Value = 40.90932046603844
$ws.Cells.Item(48, 4).Value = -0.5906795339615627
$ws.Cells.Item(48, 5).Value = 0.3489023118410489
$ws.Cells.Item(49, 3).Value = 41.283508081798
$ws.Cells.Item(49, 4).Value = -0.516491918202
$ws.Cells.Item(49, 5).Value = 0.26676390156798147
$ws.Cells.Item(50, 3).Value = 41.74291447674603
$ws.Cells.Item(50, 4).Value = -0.4570855232539728
$ws.Cells.Item(50, 5).Value = 0.2089271755683581
$ws.Cells.Item(51, 3).Value = 42.06826020706399
$ws.Cells.Item(51, 4).Value = -0.6317397929360098
$ws.Cells.Item(51, 5).Value = 0.39909516597883254

# Update TOTAL row and MSE row
$ws.Cells.Item(52, 3).Value = 1.0970331572942094
$ws.Cells.Item(52, 5).Value = 5.931608617089867
$ws.Cells.Item(53, 5).Value = 0.11863217234179733
